$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 71 (shifts old 71.. down to 72..)
$ws.Rows("71:71").Insert()

# Copy formatting from whole row 72 into new whole row 71
$ws.Rows("72:72").Copy()
$ws.Rows("71:71").PasteSpecial(-4122)  # xlPasteFormats = -4122

# Clean up any stray formatting beyond column Q to avoid bloat
$ws.Range("R71:XFD71").ClearFormats()

Write-Host "done"
